$d = $word.ActiveDocument
$before = $d.Paragraphs.Count

$endRange = $d.Content
$endRange.Collapse(0)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/></w:r><w:r><w:rPr><w:rStyle w:val="SigSignee"/></w:rPr><w:t>SIOBHIAN BROWN</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:rPr><w:rStyle w:val="Sigtitle"/></w:rPr><w:tab/><w:t>Authorised to sign by the Scottish Ministers</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/><w:rPr><w:rStyle w:val="Sigtitle"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="SigAdd"/></w:rPr><w:t>St Andrew’s House,</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SigAdd"/></w:rPr><w:t>Edinburgh</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SigDate"/></w:rPr><w:t>5th September 2024</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p>'
$endRange.InsertXML($xml)

# Indices (1-based) of the 11 newly-inserted paragraphs.
$p1  = $before + 1   # empty N3
$p2  = $before + 2   # empty N3
$p3  = $before + 3   # empty N3
$p4  = $before + 4   # empty N3
$p5  = $before + 5   # empty N3
$p6  = $before + 6   # SigBlock: SIOBHIAN BROWN
$p7  = $before + 7   # SigBlock: Authorised to sign ...
$p8  = $before + 8   # SigBlock: St Andrew's House,
$p9  = $before + 9   # SigBlock: Edinburgh
$p10 = $before + 10  # SigBlock: 5th September 2024
$p11 = $before + 11  # empty N3

# InsertXML only preserves the w:left part of w:ind (drops w:hanging), so
# restore the hanging indent on the N3 placeholder paragraphs explicitly.
foreach ($i in @($p1, $p2, $p3, $p4, $p5, $p11)) {
    $para = $d.Paragraphs.Item($i)
    $para.FirstLineIndent = -19.85
}

# InsertXML also drops w:rStyle on run properties, so reapply the character
# styles by locating each run's text within its own paragraph and setting
# Range.Style (Word maps a character-style assignment on a partial range to
# w:rPr/w:rStyle rather than the paragraph style).
function Set-RunStyle($paraIndex, $text, $styleName) {
    $para = $d.Paragraphs.Item($paraIndex)
    $r = $para.Range
    $r.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $r.Style = $styleName
}

Set-RunStyle $p6  "SIOBHIAN BROWN" "SigSignee"
Set-RunStyle $p7  "Authorised to sign by the Scottish Ministers" "Sigtitle"
Set-RunStyle $p8  "St Andrew’s House," "SigAdd"
Set-RunStyle $p9  "Edinburgh" "SigAdd"
Set-RunStyle $p10 "5th September 2024" "SigDate"
